$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Weekly Quantity": append new week row 58
# ----------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Cells.Item(58, 1).Value = 45676.99999999999
$wsWeekly.Cells.Item(58, 2).Value = 40
$wsWeekly.Cells.Item(58, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ----------------------------------------------------------------------
# Sheet "Monthly Trend": append new month row 23
# ----------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Cells.Item(23, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(23, 2).Value = 40
$wsMonthly.Cells.Item(23, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ----------------------------------------------------------------------
# Sheet "PO Forecast": new forecast model
#   - rows 2..57 keep their existing "ds" dates, Y values are replaced
#   - rows 58..65 get both new dates and new Y values
#   - row 66 is a brand-new row
# ----------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastVals = @(868, 549, 366, 876, 805, 497, 280, 252, 29, 227, 216, 38, 343, 422, 275, 374, 554, 577, 605, 585, 382, 85, 250, 694, 722, 355, 231, 484, 829, 871, 150, 89, 448, 983, 1188, 865, 646, 664, 892, 1120, 1105, 821, 558, 677, 624, 212, 256, 466, 501, 337, 886, 842, 862, 1170, 1044, 1044)

for ($i = 0; $i -lt $forecastVals.Length; $i++) {
    $wsForecast.Cells.Item($i + 2, 2).Value = $forecastVals[$i]
}

$forecastDates = @(45676.99999999999, 45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999, 45732.99999999999)
$forecastTailVals = @(934, 906, 1107, 1356, 1395, 1146, 849, 884, 1369)

for ($i = 0; $i -lt $forecastDates.Length; $i++) {
    $row = $i + 58
    $wsForecast.Cells.Item($row, 1).Value = $forecastDates[$i]
    $wsForecast.Cells.Item($row, 2).Value = $forecastTailVals[$i]
}

$wsForecast.Range("A58:A66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
